$d = $word.ActiveDocument

# Locate the end of the paragraph that currently ends the document body
# ("... Mas cada label deve se referia apenas ao seu id correto.") so the
# new content is anchored to it rather than assuming a fixed offset.
$anchor = $d.Content
$found = $anchor.Find.Execute("Mas cada label deve se referia apenas ao seu id correto.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Anchor paragraph not found"
}

$insertPos = $anchor.End
$insertRange = $d.Range($insertPos, $insertPos)

$newParagraphsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Outros inputs</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve"> – color / </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve">range / file: </w:t>
      </w:r>
      <w:r>
        <w:t>Outros tipos de inputs são: color (permite ao usuário escolher uma cor se colcoar o código hexadecimal em atributo value=”#00ff00” ele virá com essa cor de padrão). Range (alcance) é um tipo de input que permite ao usuário selecionar uma área de alcance</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> (o range padrão dele é de 0 até 100, mas podemos personalizar com atributos min=”1” e max=”5”, ai ele andará de um em um até 5).</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:tab/>
        <w:t>&lt;</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve">Input:file&gt; </w:t>
      </w:r>
      <w:r>
        <w:t>São inputs que permitem ao usuário enviar arquivos como fotos, vídeos, pedfs entre outros .jpg / .png. Como provavelmente um arquivo possui um tamanho maior do que 3.000 bites, o método de envio deve sser mudado para method=”post”.</w:t>
      </w:r>
    </w:p>
'@

$insertRange.InsertXML($newParagraphsXml)
